# "una mejora para los mensajes"
# Updates the CAADI_GROUPS Maestria roster: refreshes the header title/count,
# adds a new student row (José Ángel Calderón) with its e-mail hyperlink,
# and tidies up the leftover placeholder cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1: title / group label / student counter ----
$ws.Range("A1").Value = "Maestría"
$ws.Range("B1").Value = "ENGLISH"
$ws.Range("C1").Value = 3

# Drop the unused formatted placeholders that used to trail row 1
$ws.Range("F1").Clear()
$ws.Range("G1").Clear()

# ---- Row 2: was just formatted placeholders, now fully blank ----
$ws.Range("F2").Clear()
$ws.Range("G2").Clear()

# ---- Row 3: new student record ----
$ws.Range("A3").Value = 810069
$ws.Range("B3").Value = "José Ángel "
$ws.Range("C3").Value = "Calderón"
$ws.Range("D3").Value = "Calderón"
$ws.Range("E3").Value = "M"
$ws.Range("F3").Value = "angelcalderon2@hotmail.com"
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:angelcalderon2@hotmail.com", [Type]::Missing, [Type]::Missing, "angelcalderon2@hotmail.com")
$ws.Range("F3").NumberFormat = "General"

$ws.Range("A3:F3").VerticalAlignment = -4107
$ws.Range("B3").Font.Underline = -4142

$ws.Range("G3").Clear()

# ---- Row 4: blank templated row underneath the new record ----
$ws.Range("B4").Clear()
$ws.Range("F4").Clear()
$ws.Range("G4").Clear()

$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:placeholder@example.com", [Type]::Missing, [Type]::Missing, "")
$ws.Range("F4").Hyperlinks.Item(1).Delete()
$ws.Range("F4").ClearContents()
$ws.Range("F4").NumberFormat = "General"

$ws.Range("A4:F4").VerticalAlignment = -4107

# ---- Selection state: anchor at A1, active cell at H8 over the whole table ----
$ws.Range("A1:H8").Select()
$ws.Range("H8").Activate()

"done"
